$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bij 10.000 MB")

# Change A19, A21, A23 from the numeric value 200 to the text "sdfsd".
# This turns the shared formulas in D/E/F/G on those rows into #VALUE! errors
# since the multiplication A*B*C now includes a text operand.
$ws.Range("A19").Value = "sdfsd"
$ws.Range("A21").Value = "sdfsd"
$ws.Range("A23").Value = "sdfsd"

# Update the sheet's selection to A23 (was A12 previously).
$ws.Activate()
$ws.Range("A23").Select()
